$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 248
$ws.Range("B3").Value = 261
$ws.Range("B4").Value = 241
$ws.Range("B6").Value = 284
$ws.Range("B7").Value = 256
$ws.Range("B8").Value = 240
$ws.Range("B10").Value = 241
$ws.Range("B11").Value = 298
$ws.Range("B12").Value = 263
$ws.Range("B13").Value = 252
$ws.Range("B14").Value = 260
$ws.Range("B16").Value = 226
$ws.Range("B17").Value = 265
$ws.Range("B18").Value = 271
$ws.Range("B20").Value = 245
$ws.Range("B21").Value = 251
$ws.Range("B22").Value = 218
$ws.Range("B23").Value = 231
$ws.Range("B24").Value = 276
$ws.Range("B26").Value = 251
$ws.Range("B28").Value = 268
$ws.Range("B29").Value = 270
$ws.Range("B30").Value = 267
$ws.Range("B31").Value = 273
$ws.Range("B32").Value = 241
$ws.Range("B33").Value = 270
$ws.Range("B35").Value = 273
$ws.Range("B36").Value = 271
$ws.Range("B37").Value = 265
$ws.Range("B38").Value = 280
$ws.Range("B39").Value = 274
$ws.Range("B40").Value = 245
$ws.Range("B41").Value = 242
$ws.Range("B42").Value = 274
$ws.Range("B43").Value = 275
$ws.Range("B44").Value = 268
$ws.Range("B47").Value = 256
$ws.Range("B48").Value = 243
$ws.Range("B49").Value = 239
$ws.Range("B50").Value = 266
$ws.Range("B52").Value = 266
$ws.Range("B57").Value = 265
$ws.Range("B58").Value = 245
$ws.Range("B59").Value = 250
$ws.Range("B60").Value = 248
